$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row that only held the "5840917 - Fabricio Maciel Gomes" value
# (it had no label in column A). Everything below shifts up by one row,
# which also carries the row heights up correctly.
$ws.Rows(13).Delete()

# Fix up the cell text that the raw diff shows as re-authored content
# (these are not simple shifts - the underlying text itself changed).
$ws.Range("B10").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("C10").Value = "5840917 - Fabrício Maciel Gomes"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2021" looks like a date to Excel, so typing it directly would get
# auto-converted to a date value. Copy it from a cell that already holds it
# as plain text (row 8) so it stays a text string, matching the source data.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

$ws.Range("B18").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("C18").Value = "5840917 - Fabrício Maciel Gomes"

$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
